# Hindalco prices workbook update (2026-01-04 snapshot)
# A new "as of" row is inserted at the top of the price table (row 2),
# shifting every existing data row down by one; the previous last row
# (12-06-2025) reappears unchanged as the new last row (208).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 2, pushing all data down by one row.
$ws.Rows.Item(2).Insert()

# 2. Copy the formatting (styles) of the row that is now directly below
#    (old row 2, now row 3) onto the freshly inserted row 2, so the new
#    row looks like every other data row instead of like the header.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# 3. Populate the new top row with the latest price snapshot. The date-like
#    columns (A, E) must stay plain text (as every other row in the sheet
#    does) instead of being auto-parsed into Excel date serials, otherwise
#    "04-01-2026" (4-Jan-2026, DD-MM-YYYY) would turn into 1-Apr-2026.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "04-01-2026"
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 319.25
$ws.Range("E2").Value = "03.01.2026"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-03-january-2026.pdf"

# 4. The worksheet's Hyperlinks collection does not automatically follow
#    the row shift caused by the insert above, so rebuild it from
#    scratch based on the (already-correct) text now sitting in column F.
$ws.Hyperlinks.Delete()
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $txt = $cell.Text
    if ($txt -ne "") {
        $ws.Hyperlinks.Add($cell, $txt, "", "", $txt)
    }
}
